# Applies the "Added ABP test cases and modified IWP Bootstrap deferred
# test cases" change: refreshes the Result/Date pair on row 2 of each
# VRelayPaymentsACH_27 sheet, flipping several IWP Bootstrap deferred
# sheets (CCDeferredCorp_27, CMCAutopay*_27, PayNowNoCFCorp_27) from
# "Pass" to "Fail" as part of the re-run, and refreshing the run
# timestamp in column B for every sheet.

$wb = $excel.ActiveWorkbook

function Set-RunResult {
    param([string]$SheetName, [string]$Result, [string]$Date)

    $ws = $wb.Worksheets.Item($SheetName)
    $ws.Range("A2").Value = $Result
    $ws.Range("B2").Value = $Date
}

Set-RunResult "PayNowNoCFPC_27"   "Pass" "Sat Aug 30 00:17:12 IST 2025"
Set-RunResult "PayNowNoCFPS_27"   "Pass" "Sat Aug 30 00:17:48 IST 2025"
Set-RunResult "PayNowNoCFCorp_27" "Fail" "Sat Aug 30 00:16:36 IST 2025"
Set-RunResult "PayNowSCFPC_27"    "Pass" "Sat Aug 30 00:19:45 IST 2025"
Set-RunResult "PayNowSCFPS_27"    "Pass" "Sat Aug 30 00:20:42 IST 2025"
Set-RunResult "PayNowSCFCorp_27"  "Pass" "Sat Aug 30 00:19:10 IST 2025"
Set-RunResult "PayNowDCFPC_27"    "Pass" "Sat Aug 30 00:13:51 IST 2025"
Set-RunResult "PayNowDCFPS_27"    "Pass" "Sat Aug 30 00:14:46 IST 2025"
Set-RunResult "PayNowDCFCorp_27"  "Pass" "Sat Aug 30 00:12:59 IST 2025"
Set-RunResult "CCDeferredPS_27"   "Fail" "Sat Aug 30 00:06:16 IST 2025"
Set-RunResult "CCDeferredPC_27"   "Fail" "Sat Aug 30 00:04:03 IST 2025"
Set-RunResult "CCDeferredCorp_27" "Fail" "Sat Aug 30 00:00:35 IST 2025"
Set-RunResult "CMCAutopayPC_27"   "Fail" "Sat Aug 30 00:10:33 IST 2025"
Set-RunResult "CMCAutopayCorp_27" "Fail" "Sat Aug 30 00:09:47 IST 2025"
Set-RunResult "CMCAutopayPS_27"   "Fail" "Sat Aug 30 00:11:23 IST 2025"
